$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.231.29"
$ws.Range("E2").Value = "  +5.09%  "
$ws.Range("D3").Value = "2.357.60"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D5").Value = "'108.76"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").Value = "'309.29"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'41.13"
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("D11").Value = "'0.0917"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "'8.44"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'0.983"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "2.723.22"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'15.33"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "2.362.36"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "45.240.55"
$ws.Range("E18").Value = "  +5.89%  "
$ws.Range("D19").Value = "'7.32"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").Value = "'13.89"
$ws.Range("E20").Value = "  +6.26%  "
$ws.Range("D21").Value = "'0.0000106"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'73.23"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'3.48"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'259.02"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'11.12"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "'7.31"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("D29").Value = "'2.36"
$ws.Range("E29").Value = "  +4.70%  "
$ws.Range("D30").Value = "'0.0965"
$ws.Range("E30").Value = "  +10.21%  "
$ws.Range("D31").Value = "'22.28"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "'37.72"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").Value = "'168.95"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  +6.60%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +4.14%  "
$ws.Range("D37").Value = "'4.81"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.97"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.92"
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("D40").Value = "'0.0355"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").Value = "'1.73"
$ws.Range("E41").Value = "  +7.31%  "
$ws.Range("D42").Value = "'99.16"
$ws.Range("E42").Value = "  -4.04%  "
$ws.Range("D43").Value = "'0.232"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "'69.66"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "'12.85"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'82.10"
$ws.Range("E47").Value = "  +6.61%  "
$ws.Range("D48").Value = "'112.02"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "'5.51"
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("D50").Value = "1.675.52"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "'9.12"
$ws.Range("E51").Value = "  +3.72%  "
